# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - Updated case counters for several countries (Austria, Israel, Hungria,
#    Bulgaria, Taiwan)
#  - Suazilandia's case count grew enough to move it up in the
#    (descending, by "Casos totales") sort order, now sitting right after
#    Zimbabue and ahead of Curazao/Botsuana/Belice/... which each shift
#    down by one row
#  - Refreshed the "Datos actualizados" timestamp string in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Austria (row 19) ---
$ws.Cells.Item(19, 2).Value = 13810   # Casos totales
$ws.Cells.Item(19, 3).Value = 4       # Nuevos casos
$ws.Cells.Item(19, 5).Value = 6869    # Recuperados

# --- Israel (row 21) ---
$ws.Cells.Item(21, 2).Value = 10878   # Casos totales
$ws.Cells.Item(21, 3).Value = 135     # Nuevos casos
$ws.Cells.Item(21, 4).Value = 1388    # Casos activos
$ws.Cells.Item(21, 5).Value = 9387    # Recuperados
$ws.Cells.Item(21, 6).Value = 174     # Casos criticos
$ws.Cells.Item(21, 7).Value = 2       # Muertes hoy
$ws.Cells.Item(21, 8).Value = 103     # Muertes

# --- Hungria (row 63) ---
$ws.Cells.Item(63, 5).Value = 1193    # Recuperados
$ws.Cells.Item(63, 6).Value = 58      # Casos criticos
$ws.Cells.Item(63, 7).Value = 14      # Muertes hoy
$ws.Cells.Item(63, 8).Value = 99      # Muertes

# --- Bulgaria (row 82) ---
$ws.Cells.Item(82, 6).Value = 35      # Casos criticos

# --- Taiwan (row 100) ---
$ws.Cells.Item(100, 2).Value = 388    # Casos totales
$ws.Cells.Item(100, 3).Value = 3      # Nuevos casos
$ws.Cells.Item(100, 4).Value = 109    # Casos activos
$ws.Cells.Item(100, 5).Value = 273    # Recuperados

# --- Rows 182-189: Suazilandia's updated numbers push it above Curazao,
#     Botsuana, Belice, San Cristobal y Nieves, San Vicente y las
#     Granadinas, Nepal and Malaui, each of which shifts down one row ---

# Row 182 -> Suazilandia (new data)
$ws.Cells.Item(182, 1).Value = "Suazilandia"
$ws.Cells.Item(182, 2).Value = 14
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(182, 4).Value = 7
$ws.Cells.Item(182, 5).Value = 7
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

# Row 183 -> Curazao (previous row-182 data)
$ws.Cells.Item(183, 1).Value = "Curazao"
$ws.Cells.Item(183, 2).Value = 14
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 7
$ws.Cells.Item(183, 5).Value = 6
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 1

# Row 184 -> Botsuana (previous row-183 data)
$ws.Cells.Item(184, 1).Value = "Botsuana"
$ws.Cells.Item(184, 2).Value = 13
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 12
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 1

# Row 185 -> Belice (previous row-184 data)
$ws.Cells.Item(185, 1).Value = "Belice"
$ws.Cells.Item(185, 2).Value = 13
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 11
$ws.Cells.Item(185, 6).Value = 1
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 2

# Row 186 -> San Cristobal y Nieves (previous row-185 data)
$ws.Cells.Item(186, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(186, 2).Value = 12
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 12
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 187 -> San Vicente y las Granadinas (previous row-186 data)
$ws.Cells.Item(187, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(187, 2).Value = 12
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 1
$ws.Cells.Item(187, 5).Value = 11
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

# Row 188 -> Nepal (previous row-187 data)
$ws.Cells.Item(188, 1).Value = "Nepal"
$ws.Cells.Item(188, 2).Value = 12
$ws.Cells.Item(188, 3).Value = 3
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 5).Value = 11
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

# Row 189 -> Malaui (previous row-188 data)
$ws.Cells.Item(189, 1).Value = "Malaui"
$ws.Cells.Item(189, 2).Value = 12
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 1
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2

# --- Refresh the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 08:22"
